# New Test Case added: TCR006 - Verify registering an account by entering
# different passwords into 'Password' and 'Password Confirm' fields.
#
# This appends a new row (row 7) of data to the "Register" test-case sheet,
# reusing the formatting of the row directly above it (row 4, which carries
# the same "full data row" styling) and then writes the new, not-yet-seen
# text values (which Excel will append to the shared-string table) plus a
# couple of values that already exist elsewhere in the sheet (and so get
# reused/deduped automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")
$ws.Activate()

# --- formatting -----------------------------------------------------------
# Copy the per-column formatting of row 4 (A:H) down onto row 7 so the new
# test case matches the visual style of the existing rows. (Column D already
# carries the correct style by default, so it's left untouched.)
$cols = @("A","B","C","E","F","G","H")
foreach ($col in $cols) {
  $ws.Range($col + "4").Copy()
  $ws.Range($col + "7").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Give the new row enough height to show the wrapped test-steps text.
$ws.Rows.Item(7).RowHeight = 202.8

# --- data -------------------------------------------------------------------
# Values are written in the order the brand-new strings need to land in the
# shared-string table (Test Case ID, Test Steps, Expected Result, then the
# Test Case Description), matching how the row was originally authored.
# The other cells (B/C/E/G) reuse text already present elsewhere in the
# workbook, so Excel just references the existing shared-string entries.
$ws.Range("A7").Value = "TCR006"
$ws.Range("F7").Value = "1. Click on 'My Account' dropdown menu
2. Click on 'Register' option 
3. Enter new account details into all the fields (First Name, Last Name, E-Mail, Telephone, Password, Password Confirm, Newsletter and  Privacy Policy fields)
4. Enter any password say '12345' into the 'Password' field
5. Enter any different password say 'abcde' into the 'Passsword Confirm' field
6. Click on 'Continue' button (ER-1)"
$ws.Range("H7").Value = "1.Account should not be created, instead a warning message - 'Password confirmation does not match password!' should be displayed under 'Password Confirm' field"
$ws.Range("D7").Value = "Verify registering an account by entering different passwords into 'Password' and 'Password Confirm' fields"
$ws.Range("B7").Value = "TS001"
$ws.Range("C7").Value = "Register"
$ws.Range("E7").Value = "1. Application (https://tutorialsninja.com/demo) is opened"
$ws.Range("G7").Value = "Not Applicable"

# --- selection / view --------------------------------------------------------
# Move the cursor onto the new row, mirroring the author's edit.
$ws.Range("E7").Select()
